$d = $word.ActiveDocument

# Replace "REALIZACIÓN DEL PROTOCOLO" with "REALIZACIÓN DEL CONSENTIMIENTO"
# in the "FECHA DE REALIZACIÓN DEL PROTOCOLO:" paragraph.
$d.Content.Find.Execute("REALIZACIÓN DEL PROTOCOLO", $true, $false, $false, $false, $false,
                         $true, 1, $false, "REALIZACIÓN DEL CONSENTIMIENTO", 2)
